$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style used by the other "owner = Дечо" rows (e.g. row 8) onto row 11
$ws.Range("B8:D8").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122) # xlPasteFormats

# Set the Owner value in D11 to match the shared-string "Дечо"
$ws.Range("D11").Value = "Дечо"

# Move the selection to D11
$ws.Range("D11").Select()
